$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheet1 -> "sheet1 1"
$ws1.Name = "sheet1 1"

# Replace '..' -> '-' in the affected labels. These shared strings are used
# on both sheets (House price ranges '1..4'/'5..8' on sheet1 A22/A23 and
# sheet2 B25/B26; MortgageProductType range '3..8' on sheet1 A36 and
# sheet2 B39), so update every occurrence on both sheets.
$ws1.Range("A22").Value = "1-4"
$ws1.Range("A23").Value = "5-8"
$ws1.Range("A36").Value = "3-8"

$ws2.Range("B25").Value = "1-4"
$ws2.Range("B26").Value = "5-8"
$ws2.Range("B39").Value = "3-8"

# Leave the selection/cursor where the edit happened on each sheet, and keep
# sheet2 the active tab (matches the saved view state).
[void]$ws1.Range("B36").Select()
[void]$ws2.Range("B39").Select()
[void]$ws2.Activate()
